$d = $word.ActiveDocument

# 1. Replace the thesis/coursework title text.
$d.Content.Find.Execute(
    "Разработка параллельной программы для симуляции физической модели",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Разработка эффективного метода реализации рендеринга", 2
) | Out-Null

# 2. Remove the stray "_GoBack" bookmark left over from editing in Word.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
